# Enhance shopkeeper management: add brand field to shopkeeper queries;
# implement product fetching by brand and add shopkeeper payment details retrieval
#
# Concretely (for this Order_Items sheet), this:
#  - marks order item F7 (row 7) as deleted (IsDeleted = 1)
#  - appends 4 new order item rows (57-60)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark row 7 as deleted (IsDeleted column F)
$ws.Range("F7").Value = 1

# New order item rows to append
$newRows = @(
    @(56, 0,  "CAND355", 19, 13.73124016901582, 0),
    @(57, 48, "CAND234", 3,  2690.582959641255, 0),
    @(58, 48, "PAMP003", 7,  126, 0),
    @(59, 48, "CAND355", 19, 13.73124016901582, 0)
)

$startRow = 57
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}
